# daily auto push: 2026-02-21 22:39 UTC
# A new reading was recorded for 2026/02/22 (日) at hour 5 -> ranking 201.
# It belongs right after the existing 2026/02/22 row (row 856), so insert a
# fresh row at 857 and push every following row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 857, shifting rows 857..898 down to 858..899.
$ws.Rows.Item(857).Insert()

# Column A holds the date as literal text (e.g. "2026/12/29"), not a real
# Excel date serial. Assigning a date-shaped string straight to .Value would
# get auto-converted to a date value, so force the cell to Text first, then
# restore the style to "Normal" (no explicit style), matching the rest of
# the sheet's unstyled data cells.
$dateCell = $ws.Range("A857")
$dateCell.NumberFormat = "@"
$dateCell.Value = "2026/02/22"
$dateCell.Style = "Normal"

$ws.Range("B857").Value = "日"
$ws.Range("C857").Value = 5
$ws.Range("D857").Value = 201
